$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the header cell text directly (this also drives the table header row).
$ws.Range("B1").Value = "Nombre de campaña"
$ws.Range("F1").Value = "Ingresos"

# Keep the table (ListObject) column headers in sync with the header cells.
$table = $ws.ListObjects.Item(1)
$table.ListColumns.Item(2).Name = "Nombre de campaña"
$table.ListColumns.Item(6).Name = "Ingresos"

# Update the "Tipo de campaña" values that referenced the old CEIP text.
$ws.Range("D5").Value = "Experiencia del cliente"
$ws.Range("D8").Value = "Experiencia del cliente"
$ws.Range("D13").Value = "Experiencia del cliente"
